$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H19").Value = 2203.889
$ws.Range("I19").Value = 1916.3334
$ws.Range("J19").Value = 2347.6667
$ws.Range("K19").Value = 1916.3334
$ws.Range("L19").Value = 2347.6667
$ws.Range("M19").Value = -1741.3334
$ws.Range("N19").Value = -2697.6667

$ws.Range("H57").Value = 45000
$ws.Range("I57").Value = 30000
$ws.Range("J57").Value = 52500
$ws.Range("K57").Value = 90000
$ws.Range("L57").Value = 157500
$ws.Range("M57").Value = -89501
$ws.Range("N57").Value = -158498

$ws.Range("H86").Value = 2125.1667
$ws.Range("I86").Value = 2332.3333
$ws.Range("J86").Value = 1918
$ws.Range("K86").Value = 2332.3333
$ws.Range("L86").Value = 1918
$ws.Range("M86").Value = -1209.3333
$ws.Range("N86").Value = -4164

$ws.Range("H89").Value = 2125.1667
$ws.Range("I89").Value = 2332.3333
$ws.Range("J89").Value = 1918
$ws.Range("K89").Value = 11661.6665
$ws.Range("L89").Value = 9590
$ws.Range("M89").Value = -6045.666499999999
$ws.Range("N89").Value = -20822

$ws.Range("H113").Value = 20083.166
$ws.Range("I113").Value = 23499.8
$ws.Range("K113").Value = 23499.8
$ws.Range("M113").Value = -20245.8

$ws.Range("H125").Value = 844.125
$ws.Range("I125").Value = 815.75
$ws.Range("K125").Value = 7341.75
$ws.Range("M125").Value = -4881.75

$ws.Range("H129").Value = 912.3542
$ws.Range("J129").Value = 880.381
$ws.Range("L129").Value = 2641.143
$ws.Range("N129").Value = -12641.143

$ws.Range("H138").Value = 3403.6086
$ws.Range("J138").Value = 2037.2222
$ws.Range("L138").Value = 6111.6666
$ws.Range("N138").Value = -16391.6666

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H2").Value = 3000
$ws.Range("I2").Value = 3000
$ws.Range("J2").Value = 0
$ws.Range("K2").Value = 3000
$ws.Range("L2").Value = 0
$ws.Range("N2").Value = -2887
$ws.Range("M2").ClearContents()

$ws.Range("H32").Value = 2955.9375
$ws.Range("I32").Value = 1877.36
$ws.Range("K32").Value = 1877.36
$ws.Range("M32").Value = -1590.36

$ws.Range("H45").Value = 30001204
$ws.Range("I45").Value = 45000856
$ws.Range("K45").Value = 45000856
$ws.Range("M45").Value = -45000479

$ws.Range("H104").Value = 37187
$ws.Range("J104").Value = 37187
$ws.Range("L104").Value = 37187
$ws.Range("N104").Value = -44175

$ws.Range("H116").Value = 3000
$ws.Range("I116").Value = 3000
$ws.Range("J116").Value = 0
$ws.Range("K116").Value = 3000
$ws.Range("L116").Value = 0
$ws.Range("N116").Value = -706
$ws.Range("M116").ClearContents()

$ws.Range("H122").Value = 1869.5454
$ws.Range("I122").Value = 1923.1875
$ws.Range("J122").Value = 1726.5
$ws.Range("K122").Value = 5769.5625
$ws.Range("L122").Value = 5179.5
$ws.Range("M122").Value = -3319.5625
$ws.Range("N122").Value = -10079.5

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H3").Value = 3000
$ws.Range("I3").Value = 3000
$ws.Range("J3").Value = 0
$ws.Range("K3").Value = 3000
$ws.Range("L3").Value = 0
$ws.Range("N3").Value = -2886
$ws.Range("M3").ClearContents()

$ws.Range("H62").Value = 0
$ws.Range("I62").Value = 0
$ws.Range("K62").Value = 0
$ws.Range("M62").ClearContents()

$ws.Range("H65").Value = 0
$ws.Range("I65").Value = 0
$ws.Range("K65").Value = 0
$ws.Range("M65").ClearContents()

$ws.Range("H105").Value = 2151.56
$ws.Range("I105").Value = 2200
$ws.Range("J105").Value = 1957.8
$ws.Range("K105").Value = 2200
$ws.Range("L105").Value = 1957.8
$ws.Range("M105").Value = -453
$ws.Range("N105").Value = -5451.8

$ws.Range("H107").Value = 2125.9167
$ws.Range("I107").Value = 1728.5555
$ws.Range("K107").Value = 1728.5555
$ws.Range("M107").Value = 191.4445000000001

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H7").Value = 102.5
$ws.Range("J7").Value = 0
$ws.Range("L7").Value = 0
$ws.Range("N7").ClearContents()

$ws.Range("H31").Value = 3161.5264
$ws.Range("I31").Value = 969.9091
$ws.Range("J31").Value = 6175
$ws.Range("K31").Value = 969.9091
$ws.Range("L31").Value = 6175
$ws.Range("M31").Value = -674.9091
$ws.Range("N31").Value = -6765

$ws.Range("H33").Value = 1990
$ws.Range("I33").Value = 1990
$ws.Range("K33").Value = 1990
$ws.Range("M33").Value = -1611

$ws.Range("H34").Value = 3161.5264
$ws.Range("I34").Value = 969.9091
$ws.Range("J34").Value = 6175
$ws.Range("K34").Value = 969.9091
$ws.Range("L34").Value = 6175
$ws.Range("M34").Value = -767.9091
$ws.Range("N34").Value = -6579

$ws.Range("H122").Value = 1304.5518
$ws.Range("I122").Value = 1392.75
$ws.Range("J122").Value = 1196
$ws.Range("K122").Value = 4178.25
$ws.Range("L122").Value = 3588
$ws.Range("M122").Value = -1728.25
$ws.Range("N122").Value = -8488

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H12").Value = 109.666664
$ws.Range("J12").Value = 200
$ws.Range("L12").Value = 600
$ws.Range("N12").Value = -946

$ws.Range("H29").Value = 289.16666
$ws.Range("I29").Value = 100
$ws.Range("J29").Value = 327
$ws.Range("K29").Value = 300
$ws.Range("L29").Value = 981
$ws.Range("M29").Value = -23
$ws.Range("N29").Value = -1535

$ws.Range("H39").Value = 4259.8
$ws.Range("J39").Value = 4259.8
$ws.Range("L39").Value = 12779.4
$ws.Range("N39").Value = -13367.4

$ws.Range("H102").Value = 0
$ws.Range("J102").Value = 0
$ws.Range("N102").Value = 0
$ws.Range("L102").ClearContents()

$ws.Range("H122").Value = 1323.25
$ws.Range("J122").Value = 1347.6666
$ws.Range("L122").Value = 12128.9994
$ws.Range("N122").Value = -17028.9994

$ws.Range("H131").Value = 8076804.5
$ws.Range("J131").Value = 13542.464
$ws.Range("L131").Value = 40627.392
$ws.Range("N131").Value = -50707.392

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H21").Value = 12535004
$ws.Range("I21").Value = 25000000
$ws.Range("K21").Value = 25000000
$ws.Range("M21").Value = -24999827

$ws.Range("H29").Value = 70007.75
$ws.Range("I29").Value = 70007
$ws.Range("K29").Value = 70007
$ws.Range("M29").Value = -69717

$ws.Range("H30").Value = 12535004
$ws.Range("I30").Value = 25000000
$ws.Range("K30").Value = 25000000
$ws.Range("M30").Value = -24999895

$ws.Range("H80").Value = 2476.077
$ws.Range("I80").Value = 2025
$ws.Range("J80").Value = 2558.0908
$ws.Range("K80").Value = 2025
$ws.Range("L80").Value = 2558.0908
$ws.Range("M80").Value = -1027
$ws.Range("N80").Value = -4554.0908

$ws.Range("H83").Value = 2476.077
$ws.Range("I83").Value = 2025
$ws.Range("J83").Value = 2558.0908
$ws.Range("K83").Value = 10125
$ws.Range("L83").Value = 12790.454
$ws.Range("M83").Value = -5133
$ws.Range("N83").Value = -22774.454

$ws.Range("H102").Value = 2510.2942
$ws.Range("I102").Value = 3416.8333
$ws.Range("K102").Value = 3416.8333
$ws.Range("M102").Value = -1794.8333

$ws.Range("H122").Value = 1951.4166
$ws.Range("I122").Value = 2012.4
$ws.Range("J122").Value = 1907.8572
$ws.Range("K122").Value = 6037.200000000001
$ws.Range("L122").Value = 5723.571599999999
$ws.Range("M122").Value = -3587.200000000001
$ws.Range("N122").Value = -10623.5716

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H9").Value = 0
$ws.Range("I9").Value = 0
$ws.Range("K9").Value = 0
$ws.Range("M9").ClearContents()

$ws.Range("H23").Value = 0
$ws.Range("J23").Value = 0
$ws.Range("N23").Value = 0
$ws.Range("L23").ClearContents()

$ws.Range("H61").Value = 2762.6875
$ws.Range("I61").Value = 2373.3635
$ws.Range("J61").Value = 3619.2
$ws.Range("K61").Value = 2373.3635
$ws.Range("L61").Value = 3619.2
$ws.Range("M61").Value = -2171.3635
$ws.Range("N61").Value = -4023.2

$ws.Range("H113").Value = 2762.6875
$ws.Range("I113").Value = 2373.3635
$ws.Range("J113").Value = 3619.2
$ws.Range("K113").Value = 2373.3635
$ws.Range("L113").Value = 3619.2
$ws.Range("M113").Value = -203.3634999999999
$ws.Range("N113").Value = -7959.2

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H62").Value = 33335000
$ws.Range("I62").Value = 50001500
$ws.Range("K62").Value = 50001500
$ws.Range("M62").Value = -50000876

$ws.Range("H65").Value = 33335000
$ws.Range("I65").Value = 50001500
$ws.Range("K65").Value = 250007500
$ws.Range("M65").Value = -250004380

$ws.Range("H105").Value = 41950
$ws.Range("J105").Value = 41950
$ws.Range("L105").Value = 41950
$ws.Range("N105").Value = -48938

$ws.Range("H107").Value = 862.1818
$ws.Range("I107").Value = 721.55554
$ws.Range("J107").Value = 1495
$ws.Range("K107").Value = 2164.66662
$ws.Range("L107").Value = 4485
$ws.Range("M107").Value = -244.66662
$ws.Range("N107").Value = -8325

$ws.Range("H119").Value = 0
$ws.Range("J119").Value = 0
$ws.Range("N119").Value = 0
$ws.Range("L119").ClearContents()

$ws.Range("H122").Value = 40815.85
$ws.Range("I122").Value = 71865.45
$ws.Range("J122").Value = 2866.3333
$ws.Range("K122").Value = 215596.35
$ws.Range("L122").Value = 8598.999899999999
$ws.Range("M122").Value = -213146.35
$ws.Range("N122").Value = -13498.9999

$ws.Range("H123").Value = 47470.09
$ws.Range("J123").Value = 47470.09
$ws.Range("L123").Value = 47470.09
$ws.Range("N123").Value = -57270.09

$ws.Range("H132").Value = 2015.6842
$ws.Range("I132").Value = 1164.5454
$ws.Range("K132").Value = 3493.6362
$ws.Range("M132").Value = -963.6361999999999
